$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.110880374908447
$ws.Range("B1").Value = 4.6388258934021
$ws.Range("C1").Value = 3.481638193130493
$ws.Range("D1").Value = 0.8985042572021484
$ws.Range("E1").Value = 0.4724022746086121
